# BMZ.xlsx monthly roll-forward ("manual dislocation uploading 2021/11/06 21:00")
#
# Each month the plan/fact sheet advances: the previously-populated rows
# (2021-09) shift forward one month (+30 days) with their CarAmount figures
# zeroed out, and the next batch of previously-blank template rows
# (2021-10-31 .. 2021-11-30) gets populated with the new month's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 2-31: roll the daily schedule forward by one month (30 days) ---
# and zero out the prior month's realised CarAmount figures (column B).
$ws.Range("A2").Value = 44470
$ws.Range("B2").Value = 0
$ws.Range("A3").Value = 44471
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = 44472
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = 44473
$ws.Range("B5").Value = 0
$ws.Range("A6").Value = 44474
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = 44475
$ws.Range("B7").Value = 0
$ws.Range("A8").Value = 44476
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = 44477
$ws.Range("B9").Value = 0
$ws.Range("A10").Value = 44478
$ws.Range("B10").Value = 0
$ws.Range("A11").Value = 44479
$ws.Range("B11").Value = 0
$ws.Range("A12").Value = 44480
$ws.Range("B12").Value = 0
$ws.Range("A13").Value = 44481
$ws.Range("B13").Value = 0
$ws.Range("A14").Value = 44482
$ws.Range("B14").Value = 0
$ws.Range("A15").Value = 44483
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = 44484
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = 44485
$ws.Range("B17").Value = 0
$ws.Range("A18").Value = 44486
$ws.Range("B18").Value = 0
$ws.Range("A19").Value = 44487
$ws.Range("B19").Value = 0
$ws.Range("A20").Value = 44488
$ws.Range("B20").Value = 0
$ws.Range("A21").Value = 44489
$ws.Range("B21").Value = 0
$ws.Range("A22").Value = 44490
$ws.Range("B22").Value = 0
$ws.Range("A23").Value = 44491
$ws.Range("B23").Value = 0
$ws.Range("A24").Value = 44492
$ws.Range("B24").Value = 0
$ws.Range("A25").Value = 44493
$ws.Range("B25").Value = 0
$ws.Range("A26").Value = 44494
$ws.Range("B26").Value = 0
$ws.Range("A27").Value = 44495
$ws.Range("B27").Value = 0
$ws.Range("A28").Value = 44496
$ws.Range("B28").Value = 0
$ws.Range("A29").Value = 44497
$ws.Range("B29").Value = 0
$ws.Range("A30").Value = 44498
$ws.Range("B30").Value = 0
$ws.Range("A31").Value = 44499
$ws.Range("B31").Value = 0

# --- Rows 32-62: fill in the next month's placeholder rows with real data ---
$ws.Range("A32").Value = 44500
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = "Балхаш I"
$ws.Range("D32").Value = "Достык (эксп.)"
$ws.Range("E32").Value = "МЕДЬ"
$ws.Range("A33").Value = 44501
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = "Балхаш I"
# D33 carried a leftover cell style (s="4") from the blank template; the
# author's edit drops it back to the sheet's default (unstyled) look, same
# as every other D cell in this block, so re-stamp it from a neighbouring
# default-styled cell before writing its value.
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("D33").Value = "Достык (эксп.)"
$ws.Range("E33").Value = "МЕДЬ"
$ws.Range("A34").Value = 44502
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = "Балхаш I"
$ws.Range("D34").Value = "Достык (эксп.)"
$ws.Range("E34").Value = "МЕДЬ"
$ws.Range("A35").Value = 44503
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = "Балхаш I"
$ws.Range("D35").Value = "Достык (эксп.)"
$ws.Range("E35").Value = "МЕДЬ"
$ws.Range("A36").Value = 44504
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = "Балхаш I"
$ws.Range("D36").Value = "Достык (эксп.)"
$ws.Range("E36").Value = "МЕДЬ"
$ws.Range("A37").Value = 44505
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = "Балхаш I"
$ws.Range("D37").Value = "Достык (эксп.)"
$ws.Range("E37").Value = "МЕДЬ"
$ws.Range("A38").Value = 44506
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = "Балхаш I"
$ws.Range("D38").Value = "Достык (эксп.)"
$ws.Range("E38").Value = "МЕДЬ"
$ws.Range("A39").Value = 44507
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = "Балхаш I"
$ws.Range("D39").Value = "Достык (эксп.)"
$ws.Range("E39").Value = "МЕДЬ"
$ws.Range("A40").Value = 44508
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = "Балхаш I"
$ws.Range("D40").Value = "Достык (эксп.)"
$ws.Range("E40").Value = "МЕДЬ"
$ws.Range("A41").Value = 44509
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = "Балхаш I"
$ws.Range("D41").Value = "Достык (эксп.)"
$ws.Range("E41").Value = "МЕДЬ"
$ws.Range("A42").Value = 44510
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = "Балхаш I"
$ws.Range("D42").Value = "Достык (эксп.)"
$ws.Range("E42").Value = "МЕДЬ"
$ws.Range("A43").Value = 44511
$ws.Range("B43").Value = 33
$ws.Range("C43").Value = "Балхаш I"
$ws.Range("D43").Value = "Достык (эксп.)"
$ws.Range("E43").Value = "МЕДЬ"
$ws.Range("A44").Value = 44512
$ws.Range("B44").Value = 0
$ws.Range("C44").Value = "Балхаш I"
$ws.Range("D44").Value = "Достык (эксп.)"
$ws.Range("E44").Value = "МЕДЬ"
$ws.Range("A45").Value = 44513
$ws.Range("B45").Value = 0
$ws.Range("C45").Value = "Балхаш I"
$ws.Range("D45").Value = "Достык (эксп.)"
$ws.Range("E45").Value = "МЕДЬ"
$ws.Range("A46").Value = 44514
$ws.Range("B46").Value = 0
$ws.Range("C46").Value = "Балхаш I"
$ws.Range("D46").Value = "Достык (эксп.)"
$ws.Range("E46").Value = "МЕДЬ"
$ws.Range("A47").Value = 44515
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = "Балхаш I"
$ws.Range("D47").Value = "Достык (эксп.)"
$ws.Range("E47").Value = "МЕДЬ"
$ws.Range("A48").Value = 44516
$ws.Range("B48").Value = 30
$ws.Range("C48").Value = "Балхаш I"
$ws.Range("D48").Value = "Достык (эксп.)"
$ws.Range("E48").Value = "МЕДЬ"
$ws.Range("A49").Value = 44517
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = "Балхаш I"
$ws.Range("D49").Value = "Достык (эксп.)"
$ws.Range("E49").Value = "МЕДЬ"
$ws.Range("A50").Value = 44518
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = "Балхаш I"
$ws.Range("D50").Value = "Достык (эксп.)"
$ws.Range("E50").Value = "МЕДЬ"
$ws.Range("A51").Value = 44519
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = "Балхаш I"
$ws.Range("D51").Value = "Достык (эксп.)"
$ws.Range("E51").Value = "МЕДЬ"
$ws.Range("A52").Value = 44520
$ws.Range("B52").Value = 0
$ws.Range("C52").Value = "Балхаш I"
$ws.Range("D52").Value = "Достык (эксп.)"
$ws.Range("E52").Value = "МЕДЬ"
$ws.Range("A53").Value = 44521
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = "Балхаш I"
$ws.Range("D53").Value = "Достык (эксп.)"
$ws.Range("E53").Value = "МЕДЬ"
$ws.Range("A54").Value = 44522
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = "Балхаш I"
$ws.Range("D54").Value = "Достык (эксп.)"
$ws.Range("E54").Value = "МЕДЬ"
$ws.Range("A55").Value = 44523
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = "Балхаш I"
$ws.Range("D55").Value = "Достык (эксп.)"
$ws.Range("E55").Value = "МЕДЬ"
$ws.Range("A56").Value = 44524
$ws.Range("B56").Value = 0
$ws.Range("C56").Value = "Балхаш I"
$ws.Range("D56").Value = "Достык (эксп.)"
$ws.Range("E56").Value = "МЕДЬ"
$ws.Range("A57").Value = 44525
$ws.Range("B57").Value = 15
$ws.Range("C57").Value = "Балхаш I"
$ws.Range("D57").Value = "Достык (эксп.)"
$ws.Range("E57").Value = "МЕДЬ"
$ws.Range("A58").Value = 44526
$ws.Range("B58").Value = 0
$ws.Range("C58").Value = "Балхаш I"
$ws.Range("D58").Value = "Достык (эксп.)"
$ws.Range("E58").Value = "МЕДЬ"
$ws.Range("A59").Value = 44527
$ws.Range("B59").Value = 0
$ws.Range("C59").Value = "Балхаш I"
$ws.Range("D59").Value = "Достык (эксп.)"
$ws.Range("E59").Value = "МЕДЬ"
$ws.Range("A60").Value = 44528
$ws.Range("B60").Value = 0
$ws.Range("C60").Value = "Балхаш I"
$ws.Range("D60").Value = "Достык (эксп.)"
$ws.Range("E60").Value = "МЕДЬ"
$ws.Range("A61").Value = 44529
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = "Балхаш I"
$ws.Range("D61").Value = "Достык (эксп.)"
$ws.Range("E61").Value = "МЕДЬ"
$ws.Range("A62").Value = 44530
$ws.Range("B62").Value = 16
$ws.Range("C62").Value = "Балхаш I"
$ws.Range("D62").Value = "Достык (эксп.)"
$ws.Range("E62").Value = "МЕДЬ"

# --- Window / selection state ---
# Scroll the viewport so row 30 is the top-left visible row, and leave the
# active selection on B62 (the last newly-filled cell), matching the
# author's on-screen state after the edit.
try {
    $excel.ActiveWindow.ScrollRow = 30
} catch {
}
$ws.Range("B62").Select()
